$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-22 Thursday" "2024-02-23 Friday"

Replace-Text "725÷7=103, 4" "616÷4=154, 0"
Replace-Text "650÷6=108, 2" "656÷3=218, 2"
Replace-Text "928÷4=232, 0" "433÷8=54, 1"
Replace-Text "294÷4=73, 2" "398÷7=56, 6"
Replace-Text "644÷7=92, 0" "572÷3=190, 2"

Replace-Text "784÷8=98, 0" "486÷2=243, 0"
Replace-Text "979÷9=108, 7" "217÷5=43, 2"
Replace-Text "763÷7=109, 0" "428÷6=71, 2"
Replace-Text "732÷2=366, 0" "584÷2=292, 0"
Replace-Text "111÷4=27, 3" "131÷9=14, 5"

Replace-Text "487÷7=69, 4" "369÷8=46, 1"
Replace-Text "190÷7=27, 1" "352÷3=117, 1"
Replace-Text "322÷7=46, 0" "450÷6=75, 0"
Replace-Text "856÷4=214, 0" "968÷9=107, 5"
Replace-Text "801÷3=267, 0" "916÷9=101, 7"

Replace-Text "171÷5=34, 1" "442÷7=63, 1"
Replace-Text "923÷7=131, 6" "165÷3=55, 0"
Replace-Text "242÷2=121, 0" "531÷2=265, 1"
Replace-Text "231÷8=28, 7" "745÷2=372, 1"
Replace-Text "765÷7=109, 2" "521÷8=65, 1"

Replace-Text "537÷3=179, 0" "907÷4=226, 3"
Replace-Text "867÷7=123, 6" "718÷2=359, 0"
Replace-Text "342÷8=42, 6" "638÷5=127, 3"
Replace-Text "847÷8=105, 7" "518÷3=172, 2"
Replace-Text "695÷6=115, 5" "935÷9=103, 8"

Write-Output "Done applying replacements"
